$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy cell formatting (borders, fonts, alignment) from existing rows so the
# three new rows match the look of the rest of the table, without disturbing
# the shared styles table (PasteSpecial formats reuses existing style ids).
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

$ws.Range("A12:E12").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)

$ws.Range("A12:E12").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Fill the cell values. C14 is written before B14 (and the B15/C15 and
# B16/C16 pairs before that) so that new entries land in the shared-strings
# table in the same order as the authored workbook.
$ws.Cells.Item(14, 3).Value = "Reduce the photo size in the Home page"

$ws.Cells.Item(15, 2).Value = "Remove the status bar at the bottom from red colour on wards (Social networks, useful links, latest news...)"
$ws.Cells.Item(15, 3).Value = "Remove the status bar at the bottom from red colour on wards (Social networks, useful links, latest news...)"

$ws.Cells.Item(16, 2).Value = " Modify the courses page in the format mentioned in the attached excel sheet."
$ws.Cells.Item(16, 3).Value = " Modify the courses page in the format mentioned in the attached excel sheet."

$ws.Cells.Item(14, 2).Value = "Photo Size is large in Home page."

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(16, 1).Value = 15

$ws.Cells.Item(15, 5).Value = "Completed"

# Row heights: rows 15 and 16 wrap onto two lines like the rest of the table.
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

$ws.Range("C15").Select()
